$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 769426.7
$ws.Range("I2").Value = 769426.7
$ws.Range("K2").Value = 769426.7
$ws.Range("M2").Value = -769313.7

$ws.Range("H17").Value = 2735.4285
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2735.4285
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 8206.2855
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -8542.2855

$ws.Range("H32").Value = 1887.6471
$ws.Range("I32").Value = 1099.75
$ws.Range("K32").Value = 1099.75
$ws.Range("M32").Value = -773.75

$ws.Range("H43").Value = 5214.76
$ws.Range("I43").Value = 2717.4443
$ws.Range("J43").Value = 6619.5
$ws.Range("K43").Value = 2717.4443
$ws.Range("L43").Value = 6619.5
$ws.Range("M43").Value = -2648.4443
$ws.Range("N43").Value = -6757.5

$ws.Range("H51").Value = 4377.1113
$ws.Range("J51").Value = 4213.857
$ws.Range("L51").Value = 4213.857
$ws.Range("N51").Value = -5181.857

$ws.Range("H64").Value = 8258.031999999999
$ws.Range("I64").Value = 4285.7144
$ws.Range("J64").Value = 9416.625
$ws.Range("K64").Value = 4285.7144
$ws.Range("L64").Value = 9416.625
$ws.Range("M64").Value = -4037.7144
$ws.Range("N64").Value = -9912.625

$ws.Range("H67").Value = 8258.031999999999
$ws.Range("I67").Value = 4285.7144
$ws.Range("J67").Value = 9416.625
$ws.Range("K67").Value = 4285.7144
$ws.Range("L67").Value = 9416.625
$ws.Range("M67").Value = -3427.7144
$ws.Range("N67").Value = -11132.625

$ws.Range("H112").Value = 2615
$ws.Range("I112").Value = 1683
$ws.Range("J112").Value = 3081
$ws.Range("K112").Value = 5049
$ws.Range("L112").Value = 9243
$ws.Range("M112").Value = -3941
$ws.Range("N112").Value = -11459

$ws.Range("H116").Value = 5777.846
$ws.Range("I116").Value = 6398.4287
$ws.Range("K116").Value = 6398.4287
$ws.Range("M116").Value = -2956.4287

$ws.Range("H141").Value = 2856.0344
$ws.Range("I141").Value = 2513.04
$ws.Range("J141").Value = 4999.75
$ws.Range("K141").Value = 7539.12
$ws.Range("L141").Value = 14999.25
$ws.Range("M141").Value = -2359.12
$ws.Range("N141").Value = -25359.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1192.7778
$ws.Range("I2").Value = 963.13336
$ws.Range("J2").Value = 2341
$ws.Range("K2").Value = 963.13336
$ws.Range("L2").Value = 2341
$ws.Range("M2").Value = -850.13336
$ws.Range("N2").Value = -2567

$ws.Range("H22").Value = 7544.857
$ws.Range("I22").Value = 5562.8
$ws.Range("J22").Value = 12500
$ws.Range("K22").Value = 5562.8
$ws.Range("L22").Value = 12500
$ws.Range("M22").Value = -5263.8
$ws.Range("N22").Value = -13098

$ws.Range("H32").Value = 4603.393
$ws.Range("I32").Value = 4035.3333
$ws.Range("K32").Value = 4035.3333
$ws.Range("M32").Value = -3748.3333

$ws.Range("H97").Value = 221.625
$ws.Range("I97").Value = 169.73334
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 169.73334
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = 326.26666
$ws.Range("N97").Value = -1992

$ws.Range("H110").Value = 2099.75
$ws.Range("I110").Value = 534.6
$ws.Range("J110").Value = 4708.3335
$ws.Range("K110").Value = 534.6
$ws.Range("L110").Value = 4708.3335
$ws.Range("M110").Value = 1510.4
$ws.Range("N110").Value = -8798.333500000001

$ws.Range("H116").Value = 1192.7778
$ws.Range("I116").Value = 963.13336
$ws.Range("J116").Value = 2341
$ws.Range("K116").Value = 963.13336
$ws.Range("L116").Value = 2341
$ws.Range("M116").Value = 1330.86664
$ws.Range("N116").Value = -6929

$ws.Range("H122").Value = 2501.5334
$ws.Range("I122").Value = 2015
$ws.Range("J122").Value = 2825.889
$ws.Range("K122").Value = 6045
$ws.Range("L122").Value = 8477.667000000001
$ws.Range("M122").Value = -3595
$ws.Range("N122").Value = -13377.667

$ws.Range("H132").Value = 1922.1111
$ws.Range("I132").Value = 1884.25
$ws.Range("K132").Value = 5652.75
$ws.Range("M132").Value = -3122.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1192.7778
$ws.Range("I3").Value = 963.13336
$ws.Range("J3").Value = 2341
$ws.Range("K3").Value = 963.13336
$ws.Range("L3").Value = 2341
$ws.Range("M3").Value = -849.13336
$ws.Range("N3").Value = -2569

$ws.Range("H20").Value = 1732.909
$ws.Range("I20").Value = 1092.75
$ws.Range("J20").Value = 2098.7144
$ws.Range("K20").Value = 1092.75
$ws.Range("L20").Value = 2098.7144
$ws.Range("M20").Value = -845.75
$ws.Range("N20").Value = -2592.7144

$ws.Range("H75").Value = 11165.75
$ws.Range("I75").Value = 11165.75
$ws.Range("K75").Value = 11165.75
$ws.Range("M75").Value = -10229.75

$ws.Range("H78").Value = 11165.75
$ws.Range("I78").Value = 11165.75
$ws.Range("K78").Value = 33497.25
$ws.Range("M78").Value = -28817.25

$ws.Range("H94").Value = 2063.5
$ws.Range("I94").Value = 1459
$ws.Range("J94").Value = 2909.8
$ws.Range("K94").Value = 1459
$ws.Range("L94").Value = 2909.8
$ws.Range("M94").Value = -1008
$ws.Range("N94").Value = -3811.8

$ws.Range("H97").Value = 16049
$ws.Range("I97").Value = 9259
$ws.Range("J97").Value = 49999
$ws.Range("K97").Value = 9259
$ws.Range("L97").Value = 49999
$ws.Range("M97").Value = -8268
$ws.Range("N97").Value = -51981

$ws.Range("H122").Value = 80780
$ws.Range("J122").Value = 80780
$ws.Range("L122").Value = 80780
$ws.Range("N122").Value = -90580

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3100.5134
$ws.Range("I31").Value = 1287.1364
$ws.Range("K31").Value = 1287.1364
$ws.Range("M31").Value = -992.1364000000001

$ws.Range("H34").Value = 3100.5134
$ws.Range("I34").Value = 1287.1364
$ws.Range("K34").Value = 1287.1364
$ws.Range("M34").Value = -1085.1364

$ws.Range("H60").Value = 52497.5
$ws.Range("J60").Value = 52497.5
$ws.Range("L60").Value = 52497.5
$ws.Range("N60").Value = -53519.5

$ws.Range("H92").Value = 74996
$ws.Range("J92").Value = 74996
$ws.Range("L92").Value = 74996
$ws.Range("N92").Value = -79988

$ws.Range("H94").Value = 2545.4167
$ws.Range("I94").Value = 2119
$ws.Range("J94").Value = 2687.5557
$ws.Range("K94").Value = 2119
$ws.Range("L94").Value = 2687.5557
$ws.Range("M94").Value = -1668
$ws.Range("N94").Value = -3589.5557

$ws.Range("H99").Value = 2683.182
$ws.Range("I99").Value = 2725.25
$ws.Range("J99").Value = 2571
$ws.Range("K99").Value = 2725.25
$ws.Range("L99").Value = 2571
$ws.Range("M99").Value = -1227.25
$ws.Range("N99").Value = -5567

$ws.Range("H126").Value = 2683.182
$ws.Range("I126").Value = 2725.25
$ws.Range("J126").Value = 2571
$ws.Range("K126").Value = 8175.75
$ws.Range("L126").Value = 7713
$ws.Range("M126").Value = -5705.75
$ws.Range("N126").Value = -12653

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 2000
$ws.Range("I116").Value = 2000
$ws.Range("K116").Value = 6000
$ws.Range("M116").Value = -2558

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 4966.6665
$ws.Range("J23").Value = 3025
$ws.Range("L23").Value = 3025
$ws.Range("N23").Value = -3471

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1696.8
$ws.Range("J22").Value = 1539.6666
$ws.Range("L22").Value = 1539.6666
$ws.Range("N22").Value = -2129.6666

$ws.Range("H27").Value = 1696.8
$ws.Range("J27").Value = 1539.6666
$ws.Range("L27").Value = 1539.6666
$ws.Range("N27").Value = -1753.6666

$ws.Range("H46").Value = 735
$ws.Range("J46").Value = 488.77777
$ws.Range("L46").Value = 488.77777
$ws.Range("N46").Value = -864.7777699999999

$ws.Range("H93").Value = 5471.95
$ws.Range("I93").Value = 575
$ws.Range("K93").Value = 575
$ws.Range("M93").Value = 673

$ws.Range("H122").Value = 4212.4595
$ws.Range("I122").Value = 4493.357
$ws.Range("K122").Value = 13480.071
$ws.Range("M122").Value = -11030.071

$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 15000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1499.75
$ws.Range("I100").Value = 1333.3334
$ws.Range("J100").Value = 1999
$ws.Range("K100").Value = 2666.6668
$ws.Range("L100").Value = 3998
$ws.Range("M100").Value = -2125.6668
$ws.Range("N100").Value = -5080

$ws.Range("H107").Value = 1572.1852
$ws.Range("I107").Value = 1497.9166
$ws.Range("J107").Value = 2166.3333
$ws.Range("K107").Value = 4493.7498
$ws.Range("L107").Value = 6498.999899999999
$ws.Range("M107").Value = -2573.7498
$ws.Range("N107").Value = -10338.9999

$ws.Range("H122").Value = 3655.5
$ws.Range("I122").Value = 2320.0454
$ws.Range("J122").Value = 7328
$ws.Range("K122").Value = 6960.1362
$ws.Range("L122").Value = 21984
$ws.Range("M122").Value = -4510.1362
$ws.Range("N122").Value = -26884

$ws.Range("H132").Value = 3246.8965
$ws.Range("I132").Value = 3468
$ws.Range("K132").Value = 10404
$ws.Range("M132").Value = -7874
